$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '75.946.17'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.012.10'
$ws.Range("E3").Value = '  +3.91%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '196.95'
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '614.70'
$ws.Range("E6").Value = '  +4.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +0.69%  '
$ws.Range("E9").Value = '  +6.32%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.008.55'
$ws.Range("E10").Value = '  +3.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.437'
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.19'
$ws.Range("E13").Value = '  +7.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.554.60'
$ws.Range("E14").Value = '  +3.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.77'
$ws.Range("E15").Value = '  +4.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '75.849.41'
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000190'
$ws.Range("E17").Value = '  +2.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.016.57'
$ws.Range("E18").Value = '  +3.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.40'
$ws.Range("E19").Value = '  +2.68%  '
$ws.Range("E20").Value = '  +3.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '377.42'
$ws.Range("E21").Value = '  +3.11%  '
$ws.Range("E22").Value = '  +6.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.35'
$ws.Range("E23").Value = '  +1.90%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.165.88'
$ws.Range("E24").Value = '  +3.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.21'
$ws.Range("E25").Value = '  +0.96%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.31'
$ws.Range("E27").Value = '  +3.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.73'
$ws.Range("E28").Value = '  +2.53%  '
$ws.Range("E29").Value = '  +3.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.23'
$ws.Range("E31").Value = '  +3.70%  '
$ws.Range("E32").Value = '  +2.57%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '491.16'
$ws.Range("E33").Value = '  +0.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.91'
$ws.Range("E34").Value = '  +5.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '20.48'
$ws.Range("E36").Value = '  +3.01%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.120'
$ws.Range("E37").Value = '  +10.61%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '161.84'
$ws.Range("E38").Value = '  -2.07%  '
$ws.Range("E39").Value = '  +1.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '189.75'
$ws.Range("E40").Value = '  +7.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.376'
$ws.Range("E41").Value = '  -1.23%  '
$ws.Range("E42").Value = '  -4.66%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("E44").Value = '  +5.57%  '
$ws.Range("E45").Value = '  +18.53%  '
$ws.Range("B46").Value = 'ImmutableX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.24'
$ws.Range("E46").Value = '  +7.15%  '
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '41.07'
$ws.Range("E47").Value = '  +2.50%  '
$ws.Range("E48").Value = '  +1.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.42'
$ws.Range("E49").Value = '  +9.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.590'
$ws.Range("E50").Value = '  +3.31%  '
$ws.Range("E51").Value = '  +1.41%  '
